$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (header "K") holds the strike count values; regenerate to reflect
# the new "K" calculation instead of the old Strike# values.
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 1
